$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.561089873313904
$ws.Range("B1").Value = 1.709385871887207
$ws.Range("C1").Value = 2.017876148223877
$ws.Range("D1").Value = 3.499033451080322
$ws.Range("E1").Value = 3.527324914932251
